# collisions_table.xlsx — rework the Alien/Player/Shield/... collision matrix.
# See commit message: "Found some nasty design flaws ... collision_handler
# process in HI_Datapath_Control_Unit.vhd." The spreadsheet is updated to
# reflect the corrected (if not yet implemented) collision outcomes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text value updates -----------------------------------------------
# Row 6 (Alien row): a few outcomes are corrected.
$ws.Range("G6").Value = "Gameover"        # was "HP - 1 / Destroy"
$ws.Range("J6").Value = "Hide / _"        # was "Destroy / _"
$ws.Range("L6").Value = "Gameover"        # was "_ / Gameover"

# The three "??? / ???" placeholder cells (undetermined collisions) are
# cleared out entirely now that the real outcomes are known/pending.
$ws.Range("H8").Copy() | Out-Null
$ws.Range("I8").PasteSpecial(-4122) | Out-Null   # xlPasteFormats: match the plain bordered look
$ws.Range("I8").ClearContents()

$ws.Range("H8").Copy() | Out-Null
$ws.Range("J8").PasteSpecial(-4122) | Out-Null
$ws.Range("J8").ClearContents()

$ws.Range("H8").Copy() | Out-Null
$ws.Range("J10").PasteSpecial(-4122) | Out-Null
$ws.Range("J10").ClearContents()

$excel.CutCopyMode = $false

# --- Highlight the "_ / _" style cells with the new accent color ------
# M6, M8 and M10 move from the plain fill to a distinct (green/accent6)
# highlight fill to flag them as affected by the new collision logic.
$ws.Range("M6").Interior.ThemeColor = 10   # xlThemeColorAccent6 -> theme="9"
$ws.Range("M6").Interior.TintAndShade = 0.39997558519241921

$ws.Range("M8").Interior.ThemeColor = 10
$ws.Range("M8").Interior.TintAndShade = 0.39997558519241921

$ws.Range("M10").Interior.ThemeColor = 10
$ws.Range("M10").Interior.TintAndShade = 0.39997558519241921

# --- Selection -----------------------------------------------------------
$ws.Range("J9").Select() | Out-Null
